$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 269, shifting existing rows 269:298 down to 270:299
$ws.Rows("269:269").Insert()

# Populate the newly inserted row 269 with its data (same record, but with the
# corrected date / price values per the source update)
$ws.Range("A269").Value = 4
$ws.Range("B269").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C269").Value = "Los Lagos"
$ws.Range("D269").Value = 44946
$ws.Range("E269").Value = 10
$ws.Range("F269").Value = 100112039
$ws.Range("G269").Value = "Ciboulette"
$ws.Range("H269").Value = "Sin especificar"
$ws.Range("I269").Value = "Primera"
$ws.Range("J269").Value = 240
$ws.Range("K269").Value = 3000
$ws.Range("L269").Value = 3500
$ws.Range("M269").Value = 3250
$ws.Range("N269").Value = '$/docena de atados'
$ws.Range("O269").Value = "Región Metropolitana"
$ws.Range("P269").Value = 1083
$ws.Range("Q269").Value = 3
$ws.Range("R269").Value = "Hortaliza"
